$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.752.25'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.708.74'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '671.96'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.59'
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.09'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.445'
$ws.Range('E11').Value = '  +2.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000236'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.88'
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.718.02'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.694.22'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '16.28'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.51'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '474.67'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.82'
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.655'
$ws.Range('E21').Value = '  +0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '80.56'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.849.07'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000127'
$ws.Range('E24').Value = '  +3.74%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.96'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.12'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.69'
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.62'
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('E32').Value = '  +4.80%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.94'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.689.91'
$ws.Range('E35').Value = '  +1.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.57'
$ws.Range('E36').Value = '  +4.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.10'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.25'
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0912'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '172.62'
$ws.Range('E42').Value = '  +3.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.942'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '47.08'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('E45').Value = '  +3.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.000281'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.76'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.29'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.09'
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.268'
$ws.Range('E51').Value = '  +1.00%  '
